# Lab7-Timings.xlsx update
# Adds the "Time to Sort Size of 2 Million" column (G) data, fills in
# previously-missing Size/Doubled Size/4 second size/Doubled Time values
# for the Insertion Sort, Selection Sort and Merge Sort rows, and restores
# the previous cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - STL::sort: add "Time to Sort Size of 2 Million"
$ws.Range("G5").Value = 2.4725299999999999

# Row 6 - Bubble Sort: add "Time to Sort Size of 2 Million"
$ws.Range("G6").Value = 2.3329599999999999

# Row 7 - Quick Sort: add "Time to Sort Size of 2 Million" (formula)
$ws.Range("G7").Formula = "=((2000000^2)/(B7*B7))*4"

# Row 8 - Insertion Sort: fill in 4 second size, Doubled Time, and the new column
$ws.Range("D8").Value = 13.4815
$ws.Range("F8").Value = 51.1222
$ws.Range("G8").Formula = "=((2000000^2)/(B8*B8))*4"

# Row 9 - Selection Sort: fill in Size, 4 second size, Doubled Time, and the new column
$ws.Range("B9").Value = 12335
$ws.Range("D9").Value = 15.9385
$ws.Range("F9").Value = 63.5855
$ws.Range("G9").Formula = "=((2000000^2)/(B9*B9))*4"

# Row 10 - Merge Sort: fill in Size, 4 second size, Doubled Time, and the new column
$ws.Range("B10").Value = 700000
$ws.Range("D10").Value = 8.20656
$ws.Range("F10").Value = 16.6195
$ws.Range("G10").Value = 11.8114

# Restore the selected cell as saved in the workbook
$ws.Range("G12").Select()
